# The "Hα" emission line is now only present as an excerpt/image-backed
# note for M1-13 column D (wavelength 6584 row), so the combined label
# "[O III], Hα" in F10 is no longer needed — it becomes just "[O III]",
# matching the value already used in D10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = "[O III]"

# Reflects the cursor/selection position left by the author after the edit.
$ws.Range("G10").Select()
